# Fixed random words selection, decreased block size.
# Clears the 10 "fake/lie" placeholder strings (שקר1..שקר10) that were
# used to fill the random-word distractor columns M:W for rows 15-24,
# leaving the cells blank (style preserved) instead of holding that text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M15:W24").ClearContents()

# Update the selected/active cell & scroll position on the sheet, matching
# the new state captured after the edit (selection moved, no frozen/offset
# top-left column anymore).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 85
$ws.Range("I36").Select()
